$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of exact "Recorded By" values to their corrected form
# (moves "System" from trailing position to just after the first entry)
$map = @{
    "dnasr281@gmail.com, System" = "System, dnasr281@gmail.com";
    "backup@backdoor.com, System" = "System, backup@backdoor.com";
    "system, backup@backdoor.com, System" = "system, System, backup@backdoor.com";
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($map.ContainsKey($val)) {
        $cell.Value2 = $map[$val]
    }
}
